$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.049.25"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "2.356.15"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'544.45"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").Value = "'134.73"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("E8").Value = "  +4.13%  "
$ws.Range("D10").Value = "'5.57"
$ws.Range("E10").Value = "  +3.45%  "
$ws.Range("E11").Value = "  -1.38%  "
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("D13").Value = "'23.93"
$ws.Range("E13").Value = "  +0.95%  "
$ws.Range("D14").Value = "2.774.37"
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("D15").Value = "57.986.05"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("D17").Value = "2.343.75"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "'10.82"
$ws.Range("E18").Value = "  +3.22%  "
$ws.Range("D19").Value = "'332.47"
$ws.Range("E19").Value = "  -1.37%  "
$ws.Range("D20").Value = "'4.30"
$ws.Range("E20").Value = "  +2.03%  "
$ws.Range("D21").Value = "'6.73"
$ws.Range("E21").Value = "  -0.41%  "
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").Value = "'62.63"
$ws.Range("E23").Value = "  +1.93%  "
$ws.Range("E24").Value = "  -1.62%  "
$ws.Range("D25").Value = "'8.50"
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("D27").Value = "'1.37"
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").Value = "'171.02"
$ws.Range("E29").Value = "  -1.89%  "
$ws.Range("D30").Value = "0.0₃0741"
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").Value = "'6.14"
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").Value = "'1.03"
$ws.Range("E33").Value = "  +2.48%  "
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'4.25"
$ws.Range("E35").Value = "  +3.60%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.94%  "
$ws.Range("E37").Value = "  -1.53%  "
$ws.Range("D38").Value = "'1.62"
$ws.Range("E38").Value = "  +1.26%  "
$ws.Range("D39").Value = "'39.43"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").Value = "'144.40"
$ws.Range("E40").Value = "  -3.29%  "
$ws.Range("D41").Value = "'291.59"
$ws.Range("E41").Value = "  +1.29%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'3.67"
$ws.Range("E42").Value = "  +1.05%  "
$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").Value = "'0.379"
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("D44").Value = "'0.0943"
$ws.Range("E44").Value = "  +1.52%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "'0.0508"
$ws.Range("E45").Value = "  +0.80%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'19.14"
$ws.Range("E46").Value = "  +1.92%  "
$ws.Range("D47").Value = "'0.566"
$ws.Range("E47").Value = "  +0.55%  "
$ws.Range("D48").Value = "'0.0223"
$ws.Range("E48").Value = "  +2.65%  "
$ws.Range("D49").Value = "'17.59"
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D51").Value = "'11.07"
$ws.Range("E51").Value = "  +1.64%  "
